# Fix admin chat status filtering and completion logic
# Reorganize the "관리자 1:1문의" row into an updated "1:1문의" section with
# additional notes about chat status/completion handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old row 18 entry ("관리자 1:1문의" / "오른쪽패널에 문의 경로 기능 추가시키기")
$ws.Range("B18:E18").ClearContents()

# Row 20 ("관리자페이지" / "체크박스 디자인 수정하기" / "회원관리, 주문관리, 상품관리") is unchanged.

# New row 22: restructured "1:1문의" entry with its own 경로 (path) column.
$ws.Range("B22").Value = "1:1문의"
$ws.Range("C22").Value = "관리자페이지"
$ws.Range("D22").Value = "오른쪽패널에 문의 경로 기능 추가시키기"

# New row 23: completion / chat-content save-delete issue note.
$ws.Range("B23").Value = "1:1문의"
$ws.Range("D23").Value = "답변완료시 관리자와 유저의 채팅내용 저장/ 삭제 문제와 "

# New row 24: user/admin status note when user restarts chat.
$ws.Range("D24").Value = "유저가 다시 채팅을 시작했을 때 유저의 상태와 관리자의 상태"

# Match the final selection shown in the saved workbook.
$ws.Range("D24").Select()
